# Updates cryptos list (D = Price, E = Volume(1h)) plus a row swap for
# Fetch.AI / VeChain (rows 43-44), matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value while forcing a Text type, so values such as
# "589.14" or "71.842.47" (European thousand-separator prices) are kept as
# literal strings instead of being auto-coerced into numbers by Excel's
# type inference (which would also lose formatting like trailing zeros).
# NumberFormat is reset via ClearFormats() right after so no stray style
# lingers on the cell (matches the source file's un-styled data cells).
function Set-TextCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell 2 4 "71.842.47"
Set-TextCell 2 5 "  +3.59%  "

# Row 3 - Ethereum
Set-TextCell 3 4 "3.690.44"
Set-TextCell 3 5 "  +9.01%  "

# Row 4 - TetherUSD
Set-TextCell 4 5 "  +0.13%  "

# Row 5 - BNB
Set-TextCell 5 4 "589.14"
Set-TextCell 5 5 "  +1.56%  "

# Row 6 - Solana
Set-TextCell 6 4 "179.69"
Set-TextCell 6 5 "  +0.61%  "

# Row 7 - LidoStakedEther
Set-TextCell 7 4 "3.683.66"
Set-TextCell 7 5 "  +8.92%  "

# Row 8 - XRP
Set-TextCell 8 5 "  +5.28%  "

# Row 9 - USDC
Set-TextCell 9 5 "  +0.17%  "

# Row 10 - Dogecoin
Set-TextCell 10 4 "0.202"
Set-TextCell 10 5 "  +2.52%  "

# Row 11 - Cardano
Set-TextCell 11 5 "  +4.64%  "

# Row 12 - Avalanche
Set-TextCell 12 4 "50.05"
Set-TextCell 12 5 "  +3.57%  "

# Row 13 - ShibaInu
Set-TextCell 13 5 "  +1.51%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell 14 4 "4.288.74"
Set-TextCell 14 5 "  +9.37%  "

# Row 15 - BitcoinCash
Set-TextCell 15 4 "684.08"
Set-TextCell 15 5 "  -0.03%  "

# Row 16 - Polkadot
Set-TextCell 16 4 "8.96"
Set-TextCell 16 5 "  +4.37%  "

# Row 17 - WrappedBTC
Set-TextCell 17 4 "72.023.53"
Set-TextCell 17 5 "  +3.76%  "

# Row 18 - WrappedEther
Set-TextCell 18 4 "3.691.38"
Set-TextCell 18 5 "  +9.28%  "

# Row 19 - TRON
Set-TextCell 19 5 "  +1.95%  "

# Row 20 - Chainlink
Set-TextCell 20 4 "18.26"
Set-TextCell 20 5 "  +3.23%  "

# Row 21 - Uniswap
Set-TextCell 21 4 "11.65"
Set-TextCell 21 5 "  +3.40%  "

# Row 22 - Polygon
Set-TextCell 22 4 "0.942"
Set-TextCell 22 5 "  +3.64%  "

# Row 23 - Toncoin
Set-TextCell 23 4 "6.14"
Set-TextCell 23 5 "  +14.73%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextCell 24 4 "17.77"
Set-TextCell 24 5 "  +3.57%  "

# Row 25 - Litecoin
Set-TextCell 25 4 "103.63"
Set-TextCell 25 5 "  +2.51%  "

# Row 26 - PancakeSwap
Set-TextCell 26 4 "4.04"
Set-TextCell 26 5 "  +4.43%  "

# Row 27 - ImmutableX
Set-TextCell 27 4 "2.85"
Set-TextCell 27 5 "  +5.91%  "

# Row 28 - RenderToken
Set-TextCell 28 4 "10.20"
Set-TextCell 28 5 "  +5.14%  "

# Row 29 - EthereumClassic
Set-TextCell 29 4 "35.52"
Set-TextCell 29 5 "  +6.14%  "

# Row 30 - Filecoin
Set-TextCell 30 4 "9.20"
Set-TextCell 30 5 "  +5.49%  "

# Row 31 - NEARProtocol
Set-TextCell 31 4 "7.40"
Set-TextCell 31 5 "  +6.82%  "

# Row 32 - dogwifhat
Set-TextCell 32 5 "  +10.17%  "

# Row 33 - Bittensor
Set-TextCell 33 4 "575.85"
Set-TextCell 33 5 "  +4.92%  "

# Row 34 - Cosmos
Set-TextCell 34 4 "11.33"
Set-TextCell 34 5 "  +2.81%  "

# Row 35 - Hedera
Set-TextCell 35 5 "  +3.90%  "

# Row 36 - OKB
Set-TextCell 36 4 "59.88"
Set-TextCell 36 5 "  +3.44%  "

# Row 37 - Maker
Set-TextCell 37 4 "3.758.55"
Set-TextCell 37 5 "  +4.35%  "

# Row 38 - Dai
Set-TextCell 38 5 "  -0.02%  "

# Row 39 - Kaspa
Set-TextCell 39 5 "  +2.96%  "

# Row 40 - PEPE
Set-TextCell 40 4 "0.0₃0773"
Set-TextCell 40 5 "  +4.42%  "

# Row 41 - InjectiveProtocol
Set-TextCell 41 4 "35.54"
Set-TextCell 41 5 "  +0.53%  "

# Row 42 - Stacks
Set-TextCell 42 5 "  +4.98%  "

# Row 43 - now VeChain (was Fetch.AI)
Set-TextCell 43 2 "VeChain"
Set-TextCell 43 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 43 4 "0.0464"
Set-TextCell 43 5 "  +9.27%  "

# Row 44 - now Fetch.AI (was VeChain)
Set-TextCell 44 2 "Fetch.AI"
Set-TextCell 44 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell 44 4 "2.80"
Set-TextCell 44 5 "  +3.04%  "

# Row 45 - TheGraph
Set-TextCell 45 5 "  +4.02%  "

# Row 46 - ApeXProtocol
Set-TextCell 46 4 "3.37"
Set-TextCell 46 5 "  -0.08%  "

# Row 47 - ThetaToken
Set-TextCell 47 4 "2.86"
Set-TextCell 47 5 "  +7.43%  "

# Row 48 - Stellar
Set-TextCell 48 5 "  +4.18%  "

# Row 49 - Mantle
Set-TextCell 49 5 "  +3.88%  "

# Row 50 - FirstDigitalUSD
Set-TextCell 50 5 "  -0.14%  "

# Row 51 - Monero
Set-TextCell 51 4 "133.92"
Set-TextCell 51 5 "  +3.41%  "
